$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The locator loop used to keep writing rows for every matched item (rows 2-6).
# The fix stops after the first match, so only row 2 survives - but now
# reflects the corrected item (N305 -> "EXPRIMIDOR DE NARANJA MANUAL") and
# the updated quantity found by the (now-correct) loop.
$ws.Rows("3:6").Delete()

$ws.Range("A2").Value = "N305-EXPRIMIDOR DE NARANJA MANUAL"
$ws.Range("P2").Value = "N305-EXPRIMIDOR DE NARANJA MANUAL"
$ws.Range("G2").Value = 5
